$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "[843.   0. 640.   0. 843. 480.   0.   0.   1.]"
$ws.Range("C3").Value = "[994.   0. 640.   0. 994. 480.   0.   0.   1.]"
$ws.Range("D3").Value = "[0.05 0.   0.   0.   0.  ]"
$ws.Range("C4").Value = "[811.   0. 640.   0. 811. 480.   0.   0.   1.]"
$ws.Range("D4").Value = "[0.06 0.08 0.   0.   0.  ]"
$ws.Range("C5").Value = "[904.   0. 640.   0. 904. 480.   0.   0.   1.]"
$ws.Range("D5").Value = "[0.07 0.02 0.02 0.   0.  ]"
$ws.Range("C6").Value = "[801.   0. 640.   0. 801. 480.   0.   0.   1.]"
$ws.Range("D6").Value = "[0.07 0.07 0.05 0.06 0.07]"
$ws.Range("C7").Value = "[826.   0. 640.   0. 894. 480.   0.   0.   1.]"
$ws.Range("C8").Value = "[936.   0. 640.   0. 855. 480.   0.   0.   1.]"
$ws.Range("D8").Value = "[0.07 0.   0.   0.   0.  ]"
$ws.Range("C9").Value = "[858.   0. 640.   0. 851. 480.   0.   0.   1.]"
$ws.Range("D9").Value = "[0.09 0.   0.   0.   0.  ]"
$ws.Range("C10").Value = "[835.   0. 640.   0. 958. 480.   0.   0.   1.]"
$ws.Range("D10").Value = "[0.02 0.02 0.09 0.   0.  ]"
$ws.Range("C11").Value = "[999.   0. 640.   0. 899. 480.   0.   0.   1.]"
$ws.Range("D11").Value = "[0.02 0.01 0.05 0.08 0.06]"
$ws.Range("C12").Value = "[940.   0. 600.   0. 940. 700.   0.   0.   1.]"
$ws.Range("C13").Value = "[994.   0. 600.   0. 994. 700.   0.   0.   1.]"
$ws.Range("D13").Value = "[0.03 0.   0.   0.   0.  ]"
$ws.Range("C14").Value = "[985.   0. 600.   0. 985. 700.   0.   0.   1.]"
$ws.Range("D14").Value = "[0.07 0.07 0.   0.   0.  ]"
$ws.Range("C15").Value = "[977.   0. 600.   0. 977. 700.   0.   0.   1.]"
$ws.Range("D15").Value = "[0.06 0.07 0.01 0.   0.  ]"
$ws.Range("C16").Value = "[900.   0. 600.   0. 900. 700.   0.   0.   1.]"
$ws.Range("D16").Value = "[0.08 0.05 0.05 0.09 0.04]"
$ws.Range("C17").Value = "[963.   0. 600.   0. 959. 700.   0.   0.   1.]"
$ws.Range("C18").Value = "[898.   0. 600.   0. 931. 700.   0.   0.   1.]"
$ws.Range("D18").Value = "[0.06 0.   0.   0.   0.  ]"
$ws.Range("C19").Value = "[823.   0. 600.   0. 933. 700.   0.   0.   1.]"
$ws.Range("D19").Value = "[0.06 0.   0.   0.   0.  ]"
$ws.Range("C20").Value = "[918.   0. 600.   0. 903. 700.   0.   0.   1.]"
$ws.Range("D20").Value = "[0.01 0.06 0.04 0.   0.  ]"
$ws.Range("C21").Value = "[936.   0. 600.   0. 902. 700.   0.   0.   1.]"
$ws.Range("D21").Value = "[0.07 0.03 0.02 0.04 0.09]"
$ws.Range("C22").Value = "[921.   0. 640.   0. 923. 480.   0.   0.   1.]"
$ws.Range("C23").Value = "[984.   0. 640.   0. 898. 480.   0.   0.   1.]"
$ws.Range("D23").Value = "[0.03 0.   0.   0.  ]"
$ws.Range("C24").Value = "[927.   0. 640.   0. 819. 480.   0.   0.   1.]"
$ws.Range("D24").Value = "[0.04 0.09 0.   0.  ]"
$ws.Range("C25").Value = "[923.   0. 640.   0. 901. 480.   0.   0.   1.]"
$ws.Range("D25").Value = "[0.08 0.08 0.06 0.  ]"
$ws.Range("C26").Value = "[872.   0. 640.   0. 952. 480.   0.   0.   1.]"
$ws.Range("D26").Value = "[0.05 0.01 0.07 0.07]"
$ws.Range("C27").Value = "[8.16e+02 1.00e-02 6.40e+02 0.00e+00 8.39e+02 4.80e+02 0.00e+00 0.00e+00`n 1.00e+00]"
$ws.Range("C28").Value = "[8.63e+02 2.00e-02 6.40e+02 0.00e+00 8.54e+02 4.80e+02 0.00e+00 0.00e+00`n 1.00e+00]"
$ws.Range("D28").Value = "[0.06 0.   0.   0.  ]"
$ws.Range("C29").Value = "[9.42e+02 4.00e-02 6.40e+02 0.00e+00 9.14e+02 4.80e+02 0.00e+00 0.00e+00`n 1.00e+00]"
$ws.Range("D29").Value = "[0.02 0.03 0.   0.  ]"
$ws.Range("C30").Value = "[8.87e+02 3.00e-02 6.40e+02 0.00e+00 9.64e+02 4.80e+02 0.00e+00 0.00e+00`n 1.00e+00]"
$ws.Range("D30").Value = "[0.04 0.09 0.01 0.  ]"
$ws.Range("C31").Value = "[9.00e+02 3.00e-02 6.40e+02 0.00e+00 8.54e+02 4.80e+02 0.00e+00 0.00e+00`n 1.00e+00]"
$ws.Range("D31").Value = "[0.09 0.02 0.04 0.02]"
$ws.Range("C32").Value = "[889.   0. 600.   0. 856. 700.   0.   0.   1.]"
$ws.Range("C33").Value = "[833.   0. 600.   0. 934. 700.   0.   0.   1.]"
$ws.Range("D33").Value = "[0.03 0.   0.   0.  ]"
$ws.Range("C34").Value = "[989.   0. 600.   0. 893. 700.   0.   0.   1.]"
$ws.Range("D34").Value = "[0.07 0.06 0.   0.  ]"
$ws.Range("C35").Value = "[873.   0. 600.   0. 804. 700.   0.   0.   1.]"
$ws.Range("D35").Value = "[0.06 0.07 0.03 0.  ]"
$ws.Range("C36").Value = "[810.   0. 600.   0. 980. 700.   0.   0.   1.]"
$ws.Range("D36").Value = "[0.08 0.01 0.04 0.07]"
$ws.Range("C37").Value = "[8.68e+02 7.00e-02 6.00e+02 0.00e+00 9.76e+02 7.00e+02 0.00e+00 0.00e+00`n 1.00e+00]"
$ws.Range("C38").Value = "[9.40e+02 7.00e-02 6.00e+02 0.00e+00 8.29e+02 7.00e+02 0.00e+00 0.00e+00`n 1.00e+00]"
$ws.Range("D38").Value = "[0.05 0.   0.   0.  ]"
$ws.Range("C39").Value = "[9.52e+02 3.00e-02 6.00e+02 0.00e+00 9.13e+02 7.00e+02 0.00e+00 0.00e+00`n 1.00e+00]"
$ws.Range("D39").Value = "[0.04 0.02 0.   0.  ]"
$ws.Range("C40").Value = "[8.99e+02 4.00e-02 6.00e+02 0.00e+00 8.33e+02 7.00e+02 0.00e+00 0.00e+00`n 1.00e+00]"
$ws.Range("D40").Value = "[0.07 0.09 0.04 0.  ]"
$ws.Range("C41").Value = "[9.27e+02 3.00e-02 6.00e+02 0.00e+00 8.85e+02 7.00e+02 0.00e+00 0.00e+00`n 1.00e+00]"
$ws.Range("D41").Value = "[0.06 0.04 0.08 0.03]"
